# Averaging point of gravity and atmosphere for x and y coordinates
#
# Updates the initial-condition vector u_i (row 34, column D) that encodes the
# averaging point of the gravity/atmosphere trajectory calculation, and bumps
# the number of trajectory-difference intervals (row 35, column D).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# u_i =  [x_min,x_avg,x_max; y_min,y_avg,y_max]  -> new averaging points for x & y
$ws.Cells.Item(34, 4).Value = "[-50,0,1000000;-0,50,50]"

# t_A_int =  number of intervals
$ws.Cells.Item(35, 4).Value = 10000

# Leave the view pointed at the top of the sheet / the Name,Description header row
$ws.Range("C3").Select()
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
